# Update "want-to-go" counts (column F) across the relevant worksheets.
# This mirrors a routine data-refresh of the scraped bilibili event stats.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    "F2"  = 15
    "F3"  = 1018
    "F4"  = 13429
    "F5"  = 40
    "F6"  = 1012
    "F8"  = 1727
    "F10" = 118
    "F11" = 73
    "F14" = 13401
    "F15" = 333
    "F16" = 591
    "F17" = 8931
    "F18" = 4
    "F19" = 7994
    "F21" = 7
    "F22" = 141
    "F30" = 393
    "F31" = 201
    "F32" = 165
    "F33" = 371
}
foreach ($cell in $sheet1Updates.Keys) {
    $ws1.Range($cell).Value = $sheet1Updates[$cell]
}

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 35

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    "F3"  = 15
    "F4"  = 1018
    "F5"  = 13429
    "F6"  = 40
    "F7"  = 1012
    "F9"  = 1727
    "F11" = 118
    "F12" = 73
    "F15" = 13401
    "F16" = 333
    "F17" = 591
    "F18" = 8931
    "F19" = 4
    "F20" = 7994
    "F22" = 7
    "F23" = 141
    "F31" = 35
    "F33" = 393
    "F34" = 201
    "F35" = 165
    "F36" = 371
}
foreach ($cell in $sheet4Updates.Keys) {
    $ws4.Range($cell).Value = $sheet4Updates[$cell]
}
